$wb = $excel.ActiveWorkbook

# 1. Remove the stale _xlnm._FilterDatabase defined name (FACTORS_CCA!$A$1:$K$78)
#    left over from a filter that is no longer applied.
foreach ($n in @($wb.Names)) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.Delete()
    }
}

# 2. Update the methodology note on the EDA/disinvested-area-population comments:
#    the underlying population source moved from 2015 PBHI to 2020 Census block centroids.
$newComment = "INPUT: Percent of population living in an EDA or disinvested area, calculated from 2020 Census block centroids and EDA/disinvested layer"

$wsMuni = $wb.Worksheets.Item("FACTORS_MUNI")
$wsMuni.Range("G1").Comment.Text($newComment)

$wsCca = $wb.Worksheets.Item("FACTORS_CCA")
$wsCca.Range("D1").Comment.Text($newComment)

# 3. Refresh FACTORS_MUNI!G (pct of population in an EDA/disinvested area) with the
#    values recalculated from 2020 Census block centroids.
$wsMuni = $wb.Worksheets.Item("FACTORS_MUNI")
$wsMuni.Range("G2").Value = 0.5683715198028122
$wsMuni.Range("G3").Value = 0.044444444444444446
$wsMuni.Range("G4").Value = 0.022556785395792898
$wsMuni.Range("G6").Value = 0.17019413976002884
$wsMuni.Range("G7").Value = 0.6298035914080934
$wsMuni.Range("G10").Value = 0.030384054448225572
$wsMuni.Range("G11").Value = 0.18829826055224425
$wsMuni.Range("G13").Value = 0.09172573513930803
$wsMuni.Range("G14").Value = 0.023255813953488372
$wsMuni.Range("G16").Value = 0.8318697110011177
$wsMuni.Range("G17").Value = 0.2566310529952692
$wsMuni.Range("G19").Value = 0.7249606986899564
$wsMuni.Range("G21").Value = 0.12523456348851755
$wsMuni.Range("G23").Value = 0.2571088444576716
$wsMuni.Range("G26").Value = 0.3027544488165854
$wsMuni.Range("G27").Value = 0.675293823455864
$wsMuni.Range("G28").Value = 0.7737215033887862
$wsMuni.Range("G29").Value = 0.00430435990002777
$wsMuni.Range("G31").Value = 0.6628961581575461
$wsMuni.Range("G36").Value = 0.5444839857651246
$wsMuni.Range("G38").Value = 0.16685903547950018
$wsMuni.Range("G39").Value = 0.7264565726772504
$wsMuni.Range("G42").Value = 0.6240512993794031
$wsMuni.Range("G43").Value = 0.776783114992722
$wsMuni.Range("G44").Value = 0.5775653017390702
$wsMuni.Range("G48").Value = 0.43540983606557376
$wsMuni.Range("G50").Value = 0.000048878244293464976
$wsMuni.Range("G51").Value = 0.14049510437834842
$wsMuni.Range("G53").Value = 0.08967195609525938
$wsMuni.Range("G57").Value = 0.14597445405850845
$wsMuni.Range("G62").Value = 0.006979695431472081
$wsMuni.Range("G65").Value = 0.5655287159072101
$wsMuni.Range("G66").Value = 0.02102889186882848
$wsMuni.Range("G68").Value = 0.40736511561518696
$wsMuni.Range("G70").Value = 0.07965689412367175
$wsMuni.Range("G71").Value = 0.25718297146868574
$wsMuni.Range("G79").Value = 0.07017923864190177
$wsMuni.Range("G82").Value = 0.08774180128960687
$wsMuni.Range("G84").Value = 0.3265312273932964
$wsMuni.Range("G85").Value = 0.07270300790473258
$wsMuni.Range("G86").Value = 0.3540752712999307
$wsMuni.Range("G89").Value = 0.1706042921686747
$wsMuni.Range("G92").Value = 0.0008793069758353416
$wsMuni.Range("G93").Value = 0.4086294416243655
$wsMuni.Range("G95").Value = 0.4350947424606352
$wsMuni.Range("G96").Value = 0.8185658464462985
$wsMuni.Range("G97").Value = 1
$wsMuni.Range("G100").Value = 0.6171723210282469
$wsMuni.Range("G102").Value = 0.34036539124439846
$wsMuni.Range("G103").Value = 0.01193001060445387
$wsMuni.Range("G104").Value = 0.9856129286558928
$wsMuni.Range("G105").Value = 0.0004807692307692308
$wsMuni.Range("G108").Value = 0.09996192651818009
$wsMuni.Range("G112").Value = 0.3489698402096285
$wsMuni.Range("G120").Value = 0.3847647676939652
$wsMuni.Range("G121").Value = 0.5473809523809524
$wsMuni.Range("G126").Value = 0.33365491651205936
$wsMuni.Range("G135").Value = 0.5263103590590178
$wsMuni.Range("G144").Value = 0.009159193684371887
$wsMuni.Range("G145").Value = 0.3023653206223581
$wsMuni.Range("G148").Value = 0.5689192937043542
$wsMuni.Range("G153").Value = 0.17081738583337702
$wsMuni.Range("G154").Value = 1
$wsMuni.Range("G158").Value = 0.8008549766091305
$wsMuni.Range("G161").Value = 0.3806631762652705
$wsMuni.Range("G167").Value = 0.0743756786102063
$wsMuni.Range("G168").Value = 0.2961220698106495
$wsMuni.Range("G169").Value = 0.321044114543024
$wsMuni.Range("G170").Value = 0.1790874524714829
$wsMuni.Range("G171").Value = 0.16067272970442692
$wsMuni.Range("G172").Value = 0.07492467112515618
$wsMuni.Range("G174").Value = 0.3276721014492754
$wsMuni.Range("G176").Value = 0.01965938338535677
$wsMuni.Range("G178").Value = 0.5826912448389089
$wsMuni.Range("G180").Value = 0.01371302027142127
$wsMuni.Range("G182").Value = 0.20973520249221184
$wsMuni.Range("G184").Value = 0.07456874590581557
$wsMuni.Range("G185").Value = 0.36803056783523524
$wsMuni.Range("G187").Value = 0.5136314067611778
$wsMuni.Range("G193").Value = 0.10371548359115224
$wsMuni.Range("G194").Value = 0.1845585203510632
$wsMuni.Range("G196").Value = 0.3913653534808419
$wsMuni.Range("G198").Value = 0.7974635383639822
$wsMuni.Range("G199").Value = 0.8971734218656338
$wsMuni.Range("G200").Value = 0.017525721202340125
$wsMuni.Range("G210").Value = 0.28926391829617637
$wsMuni.Range("G212").Value = 0.3279060665362035
$wsMuni.Range("G215").Value = 0.6316434225405202
$wsMuni.Range("G221").Value = 0.0537603305785124
$wsMuni.Range("G222").Value = 0.3603592303639967
$wsMuni.Range("G224").Value = 0.35804655870445345
$wsMuni.Range("G225").Value = 0.5335185086266759
$wsMuni.Range("G226").Value = 0.63518273888155
$wsMuni.Range("G227").Value = 0.41952707856598015
$wsMuni.Range("G228").Value = 0.71015625
$wsMuni.Range("G231").Value = 0.18863610380701956
$wsMuni.Range("G234").Value = 0.09753184713375797
$wsMuni.Range("G239").Value = 0.536408106219427
$wsMuni.Range("G242").Value = 0.23069699499165275
$wsMuni.Range("G245").Value = 0.2609091138792733
$wsMuni.Range("G247").Value = 0.9984768389929218
$wsMuni.Range("G251").Value = 0.09167247324507334
$wsMuni.Range("G255").Value = 0.4489853044086774
$wsMuni.Range("G257").Value = 0.23321205587746485
$wsMuni.Range("G262").Value = 0.17963646691280885
$wsMuni.Range("G263").Value = 0.8467213757123185
$wsMuni.Range("G265").Value = 0.2565003513703443
$wsMuni.Range("G267").Value = 0.13793511721524981
$wsMuni.Range("G271").Value = 0.578225208881621
$wsMuni.Range("G278").Value = 0.4431021625652498
$wsMuni.Range("G280").Value = 0.056308878104481874
$wsMuni.Range("G282").Value = 0.1830667186890363
$wsMuni.Range("G285").Value = 0.46278645305211924

# 4. Refresh FACTORS_CCA!D (pct of population in an EDA/disinvested area) with the
#    values recalculated from 2020 Census block centroids.
$wsCca = $wb.Worksheets.Item("FACTORS_CCA")
$wsCca.Range("D2").Value = 0.8081039755351682
$wsCca.Range("D4").Value = 0.9967602591792657
$wsCca.Range("D5").Value = 0.6865540902233686
$wsCca.Range("D6").Value = 0.9674227906769464
$wsCca.Range("D7").Value = 0.9242416396532618
$wsCca.Range("D8").Value = 0.3438359061112286
$wsCca.Range("D9").Value = 0.7075323385828943
$wsCca.Range("D10").Value = 0.9434046802191612
$wsCca.Range("D11").Value = 0.005192989464223299
$wsCca.Range("D12").Value = 0.7808141950032639
$wsCca.Range("D15").Value = 0.31242359413202936
$wsCca.Range("D16").Value = 0.6920529801324503
$wsCca.Range("D18").Value = 0.29334368487721163
$wsCca.Range("D19").Value = 0.5831649499778226
$wsCca.Range("D20").Value = 0.41284446195564
$wsCca.Range("D23").Value = 0.155925820662214
$wsCca.Range("D29").Value = 0.44061627021078476
$wsCca.Range("D30").Value = 0.9041847980804425
$wsCca.Range("D31").Value = 0.9684789170982809
$wsCca.Range("D35").Value = 0.24035850081477458
$wsCca.Range("D36").Value = 0.6003273007316134
$wsCca.Range("D37").Value = 0.3764876411351846
$wsCca.Range("D38").Value = 0.7272965055450931
$wsCca.Range("D39").Value = 0.013333333333333334
$wsCca.Range("D41").Value = 0.25798883785252136
$wsCca.Range("D42").Value = 0.5925905253610549
$wsCca.Range("D43").Value = 0.9529495422357857
$wsCca.Range("D46").Value = 0.31365052393089776
$wsCca.Range("D47").Value = 0.8263903800730084
$wsCca.Range("D48").Value = 0.0677752391425944
$wsCca.Range("D50").Value = 0.44102178812922616
$wsCca.Range("D54").Value = 0.6306167777208269
$wsCca.Range("D55").Value = 0.3766284625225178
$wsCca.Range("D56").Value = 0.47652407214189896
$wsCca.Range("D58").Value = 0.48509996826404317
$wsCca.Range("D59").Value = 0.4802052785923754
$wsCca.Range("D61").Value = 0.6537175523117854
$wsCca.Range("D62").Value = 0.948603668590272
$wsCca.Range("D64").Value = 0.7800779865295995
$wsCca.Range("D68").Value = 0.25196740232940434
$wsCca.Range("D69").Value = 1
$wsCca.Range("D75").Value = 1
$wsCca.Range("D76").Value = 0.7661238038432613
$wsCca.Range("D77").Value = 0.26175368245975783

